# Auto-generated edit script applying the cryptos.xlsx price/volume update
# commit: "Updated cryptos list on Sun Sep  3 07:24:24 UTC 2023 with GitHub Actions"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "25.956.79"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.55%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.647.04"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.83%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.007"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.56%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "215.97"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.54%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5109"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.89%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.006"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.40%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2583"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.73%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06424"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.52%  "
$ws.Range("E10").Value = "  +0.48%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07792"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.29%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.323"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.653.82"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +1.26%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.5471"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.02%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0₅7901"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.15%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "64.83"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +2.16%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.022.04"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.70%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.007"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.49%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "198.85"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.55%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.471"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +3.37%  "
$ws.Range("E21").Value = "  +1.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.081"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +2.13%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.007"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.45%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.861"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.96%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "140.14"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.80%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1152"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.17%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "6.910"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.22%  "
$ws.Range("E28").Value = "  +0.71%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.242"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.31%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05030"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +0.90%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.293"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +1.12%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.205"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.02%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.547"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.68%  "
$ws.Range("E34").Value = "  -0.11%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.8955"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.58%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.593"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.138.37"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.79%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5546"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.47%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01567"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.66%  "
$ws.Range("E40").Value = "  +0.45%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "5.672"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.8179"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +1.38%  "
$ws.Range("B43").Value = "Quant"
$ws.Range("C43").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.07"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.79%  "
$ws.Range("B44").Value = "BabyDogeCoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0₈125"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +9.69%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.785.05"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +0.84%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4538"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.54%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "55.29"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +1.27%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.003"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +0.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.05094"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.09575"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.50%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.006"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.33%  "
